$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "56.895.85"
$ws.Range("E2").Value = "  +8.13%  "
$ws.Range("D3").Value = "3.237.89"
$ws.Range("E3").Value = "  +3.71%  "
$ws.Range("E4").Value = "  +0.16%  "
$c = $ws.Range("D5")
$s = $c.Style
$c.Value = "'394.73"
$c.Style = $s
$ws.Range("E5").Value = "  -0.07%  "
$c = $ws.Range("D6")
$s = $c.Style
$c.Value = "'107.12"
$c.Style = $s
$ws.Range("E6").Value = "  +3.11%  "
$ws.Range("D7").Value = "3.236.38"
$ws.Range("E7").Value = "  +3.75%  "
$c = $ws.Range("D8")
$s = $c.Style
$c.Value = "'0.566"
$c.Style = $s
$ws.Range("E9").Value = "  +0.14%  "
$c = $ws.Range("D10")
$s = $c.Style
$c.Value = "'0.617"
$c.Style = $s
$ws.Range("E10").Value = "  +2.05%  "
$c = $ws.Range("D11")
$s = $c.Style
$c.Value = "'38.84"
$c.Style = $s
$ws.Range("E11").Value = "  +1.89%  "
$c = $ws.Range("D12")
$s = $c.Style
$c.Value = "'0.0968"
$c.Style = $s
$ws.Range("E12").Value = "  +12.23%  "
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").Value = "3.761.05"
$ws.Range("E14").Value = "  +4.37%  "
$ws.Range("E15").Value = "  +3.56%  "
$c = $ws.Range("D16")
$s = $c.Style
$c.Value = "'18.98"
$c.Style = $s
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").Value = "3.245.83"
$ws.Range("E17").Value = "  +4.28%  "
$ws.Range("E18").Value = "  -2.32%  "
$c = $ws.Range("D19")
$s = $c.Style
$c.Value = "'10.61"
$c.Style = $s
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("D20").Value = "56.793.45"
$ws.Range("E20").Value = "  +8.21%  "
$c = $ws.Range("D21")
$s = $c.Style
$c.Value = "'3.30"
$c.Style = $s
$ws.Range("E21").Value = "  +1.99%  "
$c = $ws.Range("D22")
$s = $c.Style
$c.Value = "'0.0000105"
$c.Style = $s
$ws.Range("E22").Value = "  +8.29%  "
$c = $ws.Range("D23")
$s = $c.Style
$c.Value = "'12.95"
$c.Style = $s
$ws.Range("E23").Value = "  +1.33%  "
$c = $ws.Range("D24")
$s = $c.Style
$c.Value = "'298.01"
$c.Style = $s
$ws.Range("E24").Value = "  +10.68%  "
$c = $ws.Range("D25")
$s = $c.Style
$c.Value = "'73.82"
$c.Style = $s
$ws.Range("E25").Value = "  +3.90%  "
$c = $ws.Range("D26")
$s = $c.Style
$c.Value = "'3.12"
$c.Style = $s
$ws.Range("E26").Value = "  -2.77%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c = $ws.Range("D27")
$s = $c.Style
$c.Value = "'4.38"
$c.Style = $s
$ws.Range("E27").Value = "  +3.06%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D28")
$s = $c.Style
$c.Value = "'27.86"
$c.Style = $s
$ws.Range("E28").Value = "  +1.08%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D29")
$s = $c.Style
$c.Value = "'7.83"
$c.Style = $s
$ws.Range("E29").Value = "  -2.71%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D30")
$s = $c.Style
$c.Value = "'7.23"
$c.Style = $s
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D31")
$s = $c.Style
$c.Value = "'0.169"
$c.Style = $s
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("E33").Value = "  +1.78%  "
$c = $ws.Range("D34")
$s = $c.Style
$c.Value = "'10.96"
$c.Style = $s
$ws.Range("E34").Value = "  +0.56%  "
$c = $ws.Range("D35")
$s = $c.Style
$c.Value = "'37.71"
$c.Style = $s
$ws.Range("E35").Value = "  +2.59%  "
$c = $ws.Range("D36")
$s = $c.Style
$c.Value = "'0.0482"
$c.Style = $s
$ws.Range("E36").Value = "  -1.49%  "
$ws.Range("E37").Value = "  +1.62%  "
$c = $ws.Range("D38")
$s = $c.Style
$c.Value = "'51.62"
$c.Style = $s
$ws.Range("E38").Value = "  +3.05%  "
$c = $ws.Range("D39")
$s = $c.Style
$c.Value = "'0.999"
$c.Style = $s
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D40")
$s = $c.Style
$c.Value = "'3.47"
$c.Style = $s
$ws.Range("E40").Value = "  +1.64%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D41")
$s = $c.Style
$c.Value = "'3.06"
$c.Style = $s
$ws.Range("E41").Value = "  +13.88%  "
$c = $ws.Range("D42")
$s = $c.Style
$c.Value = "'133.84"
$c.Style = $s
$ws.Range("E42").Value = "  +2.93%  "
$c = $ws.Range("D43")
$s = $c.Style
$c.Value = "'1.90"
$c.Style = $s
$ws.Range("E43").Value = "  +1.54%  "
$c = $ws.Range("D44")
$s = $c.Style
$c.Value = "'0.119"
$c.Style = $s
$ws.Range("E44").Value = "  +2.88%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Range("D45")
$s = $c.Style
$c.Value = "'0.286"
$c.Style = $s
$ws.Range("E45").Value = "  -2.25%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D46")
$s = $c.Style
$c.Value = "'17.04"
$c.Style = $s
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D47")
$s = $c.Style
$c.Value = "'3.95"
$c.Style = $s
$ws.Range("E47").Value = "  -3.57%  "
$c = $ws.Range("D48")
$s = $c.Style
$c.Value = "'21.92"
$c.Style = $s
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("D49").Value = "2.134.58"
$ws.Range("E49").Value = "  +2.60%  "
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("E51").Value = "  -2.91%  "
